$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.225.93"
$ws.Range("E2").Value = "  +4.64%  "
$ws.Range("D3").Value = "2.509.83"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "495.27"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").Value = "153.82"
$ws.Range("E6").Value = "  +11.98%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("D9").Value = "2.527.95"
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +5.54%  "
$ws.Range("D12").Value = "0.338"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "2.951.60"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "57.348.74"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "21.41"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").Value = "2.517.32"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "4.59"
$ws.Range("E19").Value = "  +6.23%  "
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  +6.01%  "
$ws.Range("D21").Value = "324.44"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +5.57%  "
$ws.Range("D24").Value = "58.57"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "2.618.44"
$ws.Range("E28").Value = "  +2.93%  "
$ws.Range("D29").Value = "7.62"
$ws.Range("E29").Value = "  +4.48%  "
$ws.Range("D30").Value = "0.0₃0835"
$ws.Range("E30").Value = "  +7.80%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "151.76"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("D34").Value = "18.35"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("D36").Value = "3.85"
$ws.Range("E36").Value = "  +7.02%  "
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("D38").Value = "0.893"
$ws.Range("E38").Value = "  +5.88%  "
$ws.Range("E39").Value = "  +10.27%  "
$ws.Range("D40").Value = "34.34"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("D41").Value = "3.55"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "4.95"
$ws.Range("E45").Value = "  +6.86%  "
$ws.Range("D46").Value = "269.49"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("D47").Value = "0.0942"
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("D49").Value = "10.21"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("D51").Value = "1.903.04"
$ws.Range("E51").Value = "  -1.42%  "
